$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the style-only empty cell at E15 (author deleted its contents/formatting)
$ws.Range("E15").Clear()

# Fill in the audit data, in the same column-pass order the original author used
# (this reproduces the exact shared-string insertion order).
$ws.Range('B3').Value = 'page2.html'
$ws.Range('B5').Value = 'pas de description'
$ws.Range('B6').Value = '1 balise li inutile'
$ws.Range('B8').Value = 'toggle navigation sans icone page2'
$ws.Range('B9').Value = 'texte trop petit'
$ws.Range('B10').Value = 'balise meta keyword '
$ws.Range('B11').Value = 'html lang = default'
$ws.Range('B7').Value = 'trop de lien dans le footer (crawlabilité)'
$ws.Range('B12').Value = 'minifiez et compressez Js et css '
$ws.Range('B13').Value = 'taille image'
$ws.Range('B15').Value = 'text format image (avis)'
$ws.Range('B16').Value = 'formulaire de contact'
$ws.Range('B14').Value = 'ressource async / defer manquante sur le script Js'
$ws.Range('C3').Value = 'fait'
$ws.Range('C4').Value = 'fait '
$ws.Range('C5').Value = 'fait'
$ws.Range('C6').Value = 'fait'
$ws.Range('C7').Value = 'fait '
$ws.Range('C8').Value = 'fait '
$ws.Range('C11').Value = 'fait '
$ws.Range('D3').Value = 'contact.html'
$ws.Range('D11').Value = 'lang = fr'
$ws.Range('D5').Value = 'ajout description avec mot clef'
$ws.Range('D6').Value = 'supp'
$ws.Range('D7').Value = 'supp'
$ws.Range('B4').Value = 'pas de titre de page et titre page2'
$ws.Range('D4').Value = 'agance la panthere et Contact'
$ws.Range('D8').Value = 'supp toggle + 1partie du menu'
$ws.Range('B17').Value = 'partie responsive nul ( text et formulaire)et ajuster css '

# Widen columns B, C, D to fit the new content (pixel widths 319 / 209 / 286 at 96dpi,
# i.e. raw OOXML column widths 45.5703125 / 29.85546875 / 40.85546875).
$ws.Columns.Item(2).ColumnWidth = 44.856026785714285
$ws.Columns.Item(3).ColumnWidth = 29.141183035714285
$ws.Columns.Item(4).ColumnWidth = 40.141183035714285

# Drop the trailing blank formatted row 1000 (sheet now ends at row 999).
$ws.Rows.Item(1000).Delete()

# Leave the selection where the author left it when saving.
$ws.Range("B22").Select()
